# Updated cryptos list on Tue Jun  6 05:38:31 UTC 2023 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest
# scraped values. A handful of "Price" cells are purely numeric-looking
# strings (e.g. "0.5100", "0.000008020") that must be forced to Text via
# NumberFormat "@" before assignment, otherwise Excel auto-converts them
# to Number and silently drops the significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.780.54'
$ws.Range("E2").Value = '  -3.98%  '
$ws.Range("D3").Value = '1.817.02'
$ws.Range("E3").Value = '  -2.88%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '277.74'
$ws.Range("E5").Value = '  -7.75%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5100'
$ws.Range("E7").Value = '  -4.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3517'
$ws.Range("E8").Value = '  -6.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.56'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06658'
$ws.Range("E10").Value = '  -7.24%  '
$ws.Range("E11").Value = '  -7.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8324'
$ws.Range("E12").Value = '  -6.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07908'
$ws.Range("E13").Value = '  -3.07%  '
$ws.Range("D14").Value = '1.823.61'
$ws.Range("E15").Value = '  -3.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.62'
$ws.Range("E16").Value = '  -5.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.09'
$ws.Range("E18").Value = '  -4.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008020'
$ws.Range("E19").Value = '  -6.01%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = '25.837.17'
$ws.Range("E21").Value = '  -3.70%  '
$ws.Range("E22").Value = '  -4.84%  '
$ws.Range("E23").Value = '  -6.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.086'
$ws.Range("E24").Value = '  -4.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.21'
$ws.Range("E25").Value = '  -3.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.184'
$ws.Range("E26").Value = '  -3.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.666'
$ws.Range("E27").Value = '  -3.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.09'
$ws.Range("E28").Value = '  -5.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.57'
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.348'
$ws.Range("E30").Value = '  -8.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.239'
$ws.Range("E31").Value = '  -7.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08831'
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04862'
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7323'
$ws.Range("E34").Value = '  -8.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.131'
$ws.Range("E35").Value = '  -3.32%  '
$ws.Range("E36").Value = '  -3.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.153'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5218'
$ws.Range("E39").Value = '  -13.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.331'
$ws.Range("E40").Value = '  -10.17%  '
$ws.Range("E41").Value = '  -5.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9538'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.205'
$ws.Range("E43").Value = '  -6.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.44'
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.048'
$ws.Range("E45").Value = '  -9.13%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4600'
$ws.Range("E47").Value = '  -10.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1364'
$ws.Range("E48").Value = '  -8.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.64'
$ws.Range("E49").Value = '  -2.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.191'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.503'
$ws.Range("E51").Value = '  -7.46%  '
